$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Effects": H column values 0 -> -1 for several rows (new "no-op"
# sentinel now that the ID-like flag became a uint/-1 convention), plus
# view/selection changes.
# ---------------------------------------------------------------------------
$wsEffects = $wb.Worksheets.Item("Effects")
$wsEffects.Range("H2").Value = -1
$wsEffects.Range("H5").Value = -1
$wsEffects.Range("H6").Value = -1
$wsEffects.Range("H8").Value = -1
$wsEffects.Range("H9").Value = -1
$wsEffects.Range("H10").Value = -1
$wsEffects.Range("H11").Value = -1
$wsEffects.Range("H12").Value = -1

# ---------------------------------------------------------------------------
# Per-sheet selection / cursor updates
# ---------------------------------------------------------------------------
$wsCards = $wb.Worksheets.Item("Cards")
$wsCards.Range("B6").Select()

$wsBuffs = $wb.Worksheets.Item("Buffs")
$wsBuffs.Range("C13").Select()

$wsConditions = $wb.Worksheets.Item("Conditions")
$wsConditions.PageSetup.Orientation = 1
$wsConditions.Range("A1:G1").Select()

# Effects becomes the active / selected sheet (and its own selection moves to
# B3:H4), so select it last so it ends up as the workbook's active tab.
$wsEffects.Range("B3:H4").Select()
$wsEffects.Select()

Write-Output "done"
